$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 22728856
$ws.Range("I32").Value = 100001400
$ws.Range("J32").Value = 1637.8235
$ws.Range("K32").Value = 100001400
$ws.Range("L32").Value = 1637.8235
$ws.Range("M32").Value = -100001074
$ws.Range("N32").Value = -2289.8235
$ws.Range("H62").Value = 2477.2917
$ws.Range("I62").Value = 2616.9048
$ws.Range("J62").Value = 1500
$ws.Range("K62").Value = 2616.9048
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = -1992.9048
$ws.Range("N62").Value = -2748
$ws.Range("H65").Value = 2477.2917
$ws.Range("I65").Value = 2616.9048
$ws.Range("J65").Value = 1500
$ws.Range("K65").Value = 13084.524
$ws.Range("L65").Value = 7500
$ws.Range("M65").Value = -9964.523999999999
$ws.Range("N65").Value = -13740
$ws.Range("H98").Value = 40638.2
$ws.Range("I98").Value = 1560.6923
$ws.Range("K98").Value = 1560.6923
$ws.Range("M98").Value = -62.69229999999993
$ws.Range("H122").Value = 40638.2
$ws.Range("I122").Value = 1560.6923
$ws.Range("K122").Value = 4682.0769
$ws.Range("M122").Value = -2232.0769
$ws.Range("H129").Value = 1611.9551
$ws.Range("I129").Value = 1369.7333
$ws.Range("J129").Value = 1661.0541
$ws.Range("K129").Value = 4109.199900000001
$ws.Range("L129").Value = 4983.1623
$ws.Range("M129").Value = 890.8000999999995
$ws.Range("N129").Value = -14983.1623
$ws.Range("H130").Value = 45900
$ws.Range("J130").Value = 45900
$ws.Range("L130").Value = 45900
$ws.Range("N130").Value = -55940
$ws.Range("H137").Value = 2541.4421
$ws.Range("I137").Value = 822.17645
$ws.Range("J137").Value = 3499.7212
$ws.Range("K137").Value = 2466.52935
$ws.Range("L137").Value = 10499.1636
$ws.Range("M137").Value = 83.47064999999975
$ws.Range("N137").Value = -15599.1636

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4803.9
$ws.Range("J2").Value = 1074.75
$ws.Range("L2").Value = 1074.75
$ws.Range("N2").Value = -1300.75
$ws.Range("H80").Value = 49991
$ws.Range("J80").Value = 49991
$ws.Range("L80").Value = 49991
$ws.Range("N80").Value = -51987
$ws.Range("H83").Value = 49991
$ws.Range("J83").Value = 49991
$ws.Range("L83").Value = 149973
$ws.Range("N83").Value = -159957
$ws.Range("H116").Value = 4803.9
$ws.Range("J116").Value = 1074.75
$ws.Range("L116").Value = 1074.75
$ws.Range("N116").Value = -5662.75
$ws.Range("H132").Value = 7694135.5
$ws.Range("I132").Value = 12196357
$ws.Range("J132").Value = 2840.3333
$ws.Range("K132").Value = 36589071
$ws.Range("L132").Value = 8520.999899999999
$ws.Range("M132").Value = -36586541
$ws.Range("N132").Value = -13580.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4803.9
$ws.Range("J3").Value = 1074.75
$ws.Range("L3").Value = 1074.75
$ws.Range("N3").Value = -1302.75
$ws.Range("H9").Value = 20000
$ws.Range("J9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("N9").Value = -20336
$ws.Range("H130").Value = 48470.43
$ws.Range("J130").Value = 48470.43
$ws.Range("L130").Value = 48470.43
$ws.Range("N130").Value = -58510.43
$ws.Range("H134").Value = 3307.3103
$ws.Range("I134").Value = 2208.4138
$ws.Range("J134").Value = 4406.207
$ws.Range("K134").Value = 6625.241399999999
$ws.Range("L134").Value = 13218.621
$ws.Range("M134").Value = -4090.241399999999
$ws.Range("N134").Value = -18288.621

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1758.625
$ws.Range("I94").Value = 1500
$ws.Range("J94").Value = 1775.8667
$ws.Range("K94").Value = 1500
$ws.Range("L94").Value = 1775.8667
$ws.Range("M94").Value = -1049
$ws.Range("N94").Value = -2677.8667
$ws.Range("H99").Value = 2220.6667
$ws.Range("I99").Value = 1941.3334
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 1941.3334
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -443.3334
$ws.Range("N99").Value = -5496
$ws.Range("H122").Value = 45510788
$ws.Range("I122").Value = 62575856
$ws.Range("J122").Value = 3942.8333
$ws.Range("K122").Value = 187727568
$ws.Range("L122").Value = 11828.4999
$ws.Range("M122").Value = -187725118
$ws.Range("N122").Value = -16728.4999
$ws.Range("H126").Value = 2220.6667
$ws.Range("I126").Value = 1941.3334
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 5824.0002
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -3354.0002
$ws.Range("N126").Value = -12440
$ws.Range("H134").Value = 1548.871
$ws.Range("I134").Value = 868.16
$ws.Range("J134").Value = 4385.1665
$ws.Range("K134").Value = 2604.48
$ws.Range("L134").Value = 13155.4995
$ws.Range("M134").Value = -69.48000000000002
$ws.Range("N134").Value = -18225.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1325
$ws.Range("I25").Value = 300
$ws.Range("J25").Value = 1666.6666
$ws.Range("K25").Value = 900
$ws.Range("L25").Value = 4999.9998
$ws.Range("M25").Value = -731
$ws.Range("N25").Value = -5337.9998
$ws.Range("H30").Value = 1325
$ws.Range("I30").Value = 300
$ws.Range("J30").Value = 1666.6666
$ws.Range("K30").Value = 900
$ws.Range("L30").Value = 4999.9998
$ws.Range("M30").Value = -798
$ws.Range("N30").Value = -5203.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2577.9048
$ws.Range("I102").Value = 2651.2222
$ws.Range("K102").Value = 2651.2222
$ws.Range("M102").Value = -1029.2222
$ws.Range("H107").Value = 78414.38
$ws.Range("I107").Value = 80150.96000000001
$ws.Range("K107").Value = 80150.96000000001
$ws.Range("M107").Value = -78230.96000000001
$ws.Range("H126").Value = 7874.222
$ws.Range("J126").Value = 2141.6667
$ws.Range("L126").Value = 6425.000100000001
$ws.Range("N126").Value = -11365.0001
$ws.Range("H130").Value = 44728
$ws.Range("J130").Value = 44728
$ws.Range("L130").Value = 44728
$ws.Range("N130").Value = -54768

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1790.2
$ws.Range("I40").Value = 1811.3846
$ws.Range("J40").Value = 1652.5
$ws.Range("K40").Value = 1811.3846
$ws.Range("L40").Value = 1652.5
$ws.Range("M40").Value = -1675.3846
$ws.Range("N40").Value = -1924.5
$ws.Range("H136").Value = 1958.8158
$ws.Range("J136").Value = 4957.6665
$ws.Range("L136").Value = 14872.9995
$ws.Range("N136").Value = -19972.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1562.1666
$ws.Range("I132").Value = 1140.8823
$ws.Range("J132").Value = 2585.2856
$ws.Range("K132").Value = 3422.6469
$ws.Range("L132").Value = 7755.8568
$ws.Range("M132").Value = -892.6468999999997
$ws.Range("N132").Value = -12815.8568
$ws.Range("H136").Value = 25493.936
$ws.Range("I136").Value = 79679.28999999999
$ws.Range("J136").Value = 2506.2122
$ws.Range("K136").Value = 239037.87
$ws.Range("L136").Value = 7518.6366
$ws.Range("M136").Value = -236487.87
$ws.Range("N136").Value = -12618.6366

